$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: "I was not particularly close..." -> double em-dash variant
$ws.Range("B4").Value = "I was not particularly close to the Lord" + [char]0x2014 + [char]0x2014 + "we only met twice in total."

# Row 5: shift up - now holds what used to be row4's content (Qingliu Manor visit)
$ws.Range("B5").Value = "The first time was during my initial visit to Qingliu Manor."

# Row 6: shift up - now holds what used to be row5's content (lunch today)
$ws.Range("B6").Value = "The second was during lunch today. I haven" + [char]0x2019 + "t seen him since."

# Row 7: shift up - now holds what used to be row6's content (please tell me...)
$ws.Range("B7").Value = "Please tell me about your actions before and during the banquet."

# Row 8: shift up - now holds what used to be row7's content (after lunch...)
$ws.Range("B8").Value = "After lunch, I remained in my room, meditating, until Steward He came to inform me of the banquet" + [char]0x2019 + "s time and location."

# Row 9: new content about 6.15 PM, and row height 51 -> 34
$ws.Range("B9").Value = "At around 6.15 PM, I left my room just before the rain began and happened to run into Ming on the way."
$ws.Rows.Item(9).RowHeight = 34

# Row 10: shift up - now holds what used to be row9's content (two of us proceeded...)
$ws.Range("B10").Value = "The two of us proceeded to the banquet hall together."

# Row 11: shift up - now holds what used to be row10's content (I did not leave...)
$ws.Range("B11").Value = "I did not leave after that."

# Row 12: wrap "So, Kong was also..." text in green color tag, row height 34 -> 51
$ws.Range("B12").Value = " <color=#00CC00>(So, Kong was also someone who arrived early at the banquet and stayed the entire time.)</color>"
$ws.Rows.Item(12).RowHeight = 51

# Row 13: wrap "If that's the case..." text in green color tag (height stays 34)
$ws.Range("B13").Value = " <color=#00CC00>(If that" + [char]0x2019 + "s the case, it seems he didn" + [char]0x2019 + "t have the opportunity to commit the crime.)</color>"

# Row 14: A14/B14 values unaffected in meaning (Goto / 10_Ran) - no content change needed

# Move active-cell selection down to B15 (matches author's final cursor position)
$ws.Range("B15").Select() | Out-Null
